$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: 2025-12-02
$ws.Range("A38").Value = [DateTime]"2025-12-02"
$ws.Range("B38").Value = 607
$ws.Range("C38").Value = 21
$ws.Range("D38").Value = 586

# Row 39: 2025-12-03
$ws.Range("A39").Value = [DateTime]"2025-12-03"
$ws.Range("B39").Value = 740
$ws.Range("C39").Value = 36
$ws.Range("D39").Value = 704

# Update the selection to A39:D39
$ws.Range("A39:D39").Select()
